# Update "想去人数" (interest count) values in column F across sheets
# 展览 (ws1), 本地生活 (ws3), 全部类型 (ws4) per the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 214
$ws1.Range("F4").Value  = 602
$ws1.Range("F6").Value  = 394
$ws1.Range("F7").Value  = 532
$ws1.Range("F12").Value = 562
$ws1.Range("F14").Value = 1713
$ws1.Range("F15").Value = 298
$ws1.Range("F16").Value = 1650
$ws1.Range("F18").Value = 479
$ws1.Range("F19").Value = 30
$ws1.Range("F20").Value = 118

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5247

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5247
$ws4.Range("F6").Value  = 214
$ws4.Range("F12").Value = 602
$ws4.Range("F16").Value = 394
$ws4.Range("F17").Value = 532
$ws4.Range("F25").Value = 562
$ws4.Range("F28").Value = 1713
$ws4.Range("F29").Value = 298
$ws4.Range("F30").Value = 1650
$ws4.Range("F33").Value = 479
$ws4.Range("F34").Value = 30
$ws4.Range("F35").Value = 118
